$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'305.63"
$ws.Range("E2").Value = "'-0.66%"
$ws.Range("D3").Value = "'38.92"
$ws.Range("E3").Value = "'7.48%"
$ws.Range("D4").Value = "'5.116"
$ws.Range("E4").Value = "'1.21%"
$ws.Range("D5").Value = "'0.08078"
$ws.Range("E5").Value = "'-0.42%"
$ws.Range("D6").Value = "'1.930"
$ws.Range("E6").Value = "'-3.28%"
$ws.Range("D7").Value = "'4.197"
$ws.Range("E7").Value = "'0.77%"
$ws.Range("D8").Value = "'8.014"
$ws.Range("E8").Value = "'1.95%"
$ws.Range("D9").Value = "'0.9287"
$ws.Range("E9").Value = "'0.13%"
$ws.Range("D10").Value = "'0.1485"
$ws.Range("E10").Value = "'1.88%"
$ws.Range("D11").Value = "'0.1925"
$ws.Range("E11").Value = "'-0.28%"
$ws.Range("D12").Value = "'0.09021"
$ws.Range("E12").Value = "'-1.19%"
$ws.Range("D13").Value = "'0.03520"
$ws.Range("E13").Value = "'2.28%"
$ws.Range("D15").Value = "'0.001392"
$ws.Range("E15").Value = "'-1.68%"
$ws.Range("D16").Value = "'0.005880"
$ws.Range("E16").Value = "'-7.17%"
$ws.Range("D17").Value = "'3.781"
$ws.Range("E17").Value = "'-1.57%"
$ws.Range("D18").Value = "'3.423"
$ws.Range("E18").Value = "'-0.80%"
$ws.Range("D19").Value = "'0.3430"
$ws.Range("E19").Value = "'-0.73%"
$ws.Range("E20").Value = "'0.03%"
$ws.Range("D21").Value = "'4.685"
$ws.Range("E21").Value = "'-2.89%"
$ws.Range("D22").Value = "'0.2417"
$ws.Range("E22").Value = "'3.08%"
$ws.Range("D23").Value = "'0.04377"
$ws.Range("E23").Value = "'-0.15%"
$ws.Range("E24").Value = "'0.32%"
$ws.Range("D25").Value = "'0.004272"
$ws.Range("E25").Value = "'2.18%"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("D39").Value = "'0.02031"
$ws.Range("E39").Value = "'-0.83%"
$ws.Range("D40").Value = "'0.05053"
$ws.Range("E40").Value = "'-1.42%"
$ws.Range("D41").Value = "'0.007534"
$ws.Range("E41").Value = "'0.86%"
$ws.Range("D42").Value = "'0.009733"
$ws.Range("E42").Value = "'-3.33%"
$ws.Range("D43").Value = "'0.1345"
$ws.Range("E43").Value = "'-1.91%"
$ws.Range("D44").Value = "'0.002105"
$ws.Range("E44").Value = "'-0.85%"
$ws.Range("D45").Value = "'0.009905"
$ws.Range("E45").Value = "'0.38%"
$ws.Range("D46").Value = "'0.00006208"
$ws.Range("E46").Value = "'-1.55%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'0.002875"
$ws.Range("D49").Value = "'0.001804"
$ws.Range("E49").Value = "'12.56%"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E51").Value = "'0.06%"
